# edit.ps1
# Applies the "add scenario files for runs march 8th 2023" change to
# mansoni_coverage_scenario_3a_1.xlsx:
#   - Row 2 (5yr platform coverage): 0.6 -> 0.736 across H/J/L/N/P/R/T2, and
#     the trailing V2 cell (which held the same 0.6) is cleared/removed.
#   - Row 3 (2yr platform coverage): 0.75 -> 0.92 across X3:AZ3.
#   - Row 4 & 5 (15-50 / 50-65 yr platform coverage): 0.5 -> 0.613 across
#     X:AZ, plus a newly-styled (but empty) W column cell.
#   - A brand-new row 12 "Vector Control" product row with 0.25 coverage
#     values in every whole-year column (X, Z, AB, ... AZ).
#   - New font / cell style (size 11, black Calibri) used by the row 4/5
#     edits above.
#   - Refreshed sheet/workbook selection state (Platform Coverage tab made
#     active, MarketShare selection narrowed).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Platform Coverage")
$ws2 = $wb.Worksheets.Item("MarketShare")

# ---------------------------------------------------------------------
# Row 2: 5yr-olds platform coverage 0.6 -> 0.736 (H,J,L,N,P,R,T), and the
# last occurrence (V2) is removed entirely rather than updated.
# ---------------------------------------------------------------------
foreach ($colLetter in @("H","J","L","N","P","R","T")) {
    $ws1.Range($colLetter + "2").Value = 0.736
}
$ws1.Range("V2").ClearContents()

# ---------------------------------------------------------------------
# Row 3: 2yr-olds platform coverage 0.75 -> 0.92 across X3:AZ3.
# ---------------------------------------------------------------------
for ($c = 24; $c -le 52; $c++) {
    $ws1.Cells.Item(3, $c).Value = 0.92
}

# ---------------------------------------------------------------------
# Row 4 & 5: 0.5 -> 0.613 across X:AZ, formatted with the new font style
# (size 11, solid black RGB Calibri) -- also extends one column earlier
# to W with no value but the same new style.
# ---------------------------------------------------------------------
foreach ($r in @(4, 5)) {
    $f = $ws1.Cells.Item($r, 23).Font
    $f.Color = 0
    $f.Size = 11
    for ($c = 24; $c -le 52; $c++) {
        $cell = $ws1.Cells.Item($r, $c)
        $cell.Value = 0.613
        $cell.Font.Color = 0
        $cell.Font.Size = 11
    }
}

# ---------------------------------------------------------------------
# New row 12: "Vector Control" product, 0.25 coverage in whole-year
# columns only (X, Z, AB, AD, ... AZ).
# ---------------------------------------------------------------------
$ws1.Range("B12").Value = "Vector Control"
for ($c = 24; $c -le 52; $c += 2) {
    $ws1.Cells.Item(12, $c).Value = 0.25
}

# ---------------------------------------------------------------------
# Selection / active-sheet state: Platform Coverage becomes the active
# tab (selection AH11); MarketShare's selection narrows to D1:AV3.
# ---------------------------------------------------------------------
$ws2.Range("D1:AV3").Select()
$ws1.Activate()
$ws1.Range("AH11").Select()
